# conductor_grid.xlsx update:
# - rename DXINCRE -> DXINCRE_LEFT (row 11), update its description
# - insert a new DXINCRE_RIGHT row (row 12) with analogous description
# - insert a new MAXNOD row (row 13)
# Rows 4-10 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: DXINCRE -> DXINCRE_LEFT -------------------------------------
$ws.Range("A11").Value = "DXINCRE_LEFT"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "float"
$ws.Range("D11").Value = "size increase ratio for the spatial mesh, used for the region to the left of the refined region."
$ws.Range("E11").Value = 1.2
$ws.Range("E11").NumberFormat = "0.00E+00"
$ws.Rows(11).RowHeight = 30

# --- Row 12 (new): DXINCRE_RIGHT -----------------------------------------
$ws.Range("A12").Value = "DXINCRE_RIGHT"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "float"
$ws.Range("D12").Value = "size increase ratio for the spatial mesh, used for the region to the right of the refined region."
$ws.Range("E12").Value = 1.2
$ws.Range("E12").NumberFormat = "0.00E+00"
$ws.Rows(12).RowHeight = 30

# --- Row 13 (new): MAXNOD -------------------------------------------------
$ws.Range("A13").Value = "MAXNOD"
$ws.Range("B13").Value = "-"
$ws.Range("C13").Value = "integer"
$ws.Range("D13").Value = " maximum number of nodes for conductor spatial discretization"
$ws.Range("E13").Value = 10001

# matches the author's on-screen selection after editing the two new rows
$ws.Range("A11:XFD12").Select() | Out-Null
